$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing player table body (keep header row 1 intact)
$ws.Range("A2:C18").ClearContents()

$rows = @(
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Cason Wallace", "PG,SG", "Oklahoma City Thunder"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Royce O'Neale", "SF,PF", "Phoenix Suns"),
    @("Jordan Clarkson", "SG,SF", "Utah Jazz"),
    @("Bam Adebayo", "PF,C", "Miami Heat"),
    @("Jaylen Clark", "SG", "Minnesota Timberwolves"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Kyle Kuzma", "SF,PF", "Milwaukee Bucks"),
    @("Moses Moody", "SG,SF", "Golden State Warriors"),
    @("Zach LaVine", "SG,SF", "Sacramento Kings"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Rui Hachimura", "SF,PF", "Los Angeles Lakers")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
